$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "67.119.66"
$ws.Cells.Item(2, 5).Value = "  -0.90%  "
Set-TextValue 3 4 "2.466.32"
$ws.Cells.Item(3, 5).Value = "  -2.71%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
Set-TextValue 5 4 "583.11"
$ws.Cells.Item(5, 5).Value = "  -1.69%  "
Set-TextValue 6 4 "169.13"
$ws.Cells.Item(6, 5).Value = "  -1.93%  "
$ws.Cells.Item(7, 5).Value = "  +0.11%  "
$ws.Cells.Item(8, 5).Value = "  -2.28%  "
Set-TextValue 9 4 "2.465.47"
$ws.Cells.Item(9, 5).Value = "  -2.77%  "
$ws.Cells.Item(10, 5).Value = "  -2.89%  "
$ws.Cells.Item(11, 5).Value = "  -0.10%  "
$ws.Cells.Item(12, 5).Value = "  -2.88%  "
$ws.Cells.Item(13, 5).Value = "  -3.74%  "
$ws.Cells.Item(14, 5).Value = "  -3.42%  "
Set-TextValue 15 4 "2.866.48"
$ws.Cells.Item(15, 5).Value = "  -1.80%  "
Set-TextValue 16 4 "67.095.55"
$ws.Cells.Item(16, 5).Value = "  -0.75%  "
$ws.Cells.Item(17, 5).Value = "  -4.50%  "
Set-TextValue 18 4 "2.460.05"
$ws.Cells.Item(18, 5).Value = "  -3.41%  "
Set-TextValue 19 4 "11.14"
$ws.Cells.Item(19, 5).Value = "  -5.88%  "
$ws.Cells.Item(20, 5).Value = "  -3.47%  "
Set-TextValue 21 4 "353.71"
$ws.Cells.Item(21, 5).Value = "  -4.25%  "
$ws.Cells.Item(22, 5).Value = "  -2.84%  "
$ws.Cells.Item(23, 5).Value = "  +0.02%  "
Set-TextValue 24 4 "69.02"
$ws.Cells.Item(24, 5).Value = "  -3.88%  "
$ws.Cells.Item(25, 5).Value = "  -7.68%  "
$ws.Cells.Item(26, 5).Value = "  -7.40%  "
Set-TextValue 27 4 "9.25"
$ws.Cells.Item(27, 5).Value = "  -7.27%  "
$ws.Cells.Item(28, 5).Value = "  -1.39%  "
Set-TextValue 29 4 "2.567.78"
$ws.Cells.Item(29, 5).Value = "  -3.15%  "
Set-TextValue 30 4 "0.0₃0906"
$ws.Cells.Item(30, 5).Value = "  -6.12%  "
Set-TextValue 31 4 "516.57"
$ws.Cells.Item(31, 5).Value = "  -4.02%  "
Set-TextValue 32 4 "7.75"
$ws.Cells.Item(32, 5).Value = "  -7.79%  "
$ws.Cells.Item(33, 5).Value = "  -5.51%  "
$ws.Cells.Item(34, 5).Value = "  -6.04%  "
$ws.Cells.Item(35, 5).Value = "  +0.10%  "
$ws.Cells.Item(36, 5).Value = "  -7.04%  "
Set-TextValue 37 4 "157.58"
$ws.Cells.Item(37, 5).Value = "  -0.66%  "
$ws.Cells.Item(38, 5).Value = "  +0.26%  "
$ws.Cells.Item(39, 5).Value = "  -3.55%  "
$ws.Cells.Item(40, 5).Value = "  -5.75%  "
$ws.Cells.Item(41, 5).Value = "  +0.23%  "
$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue 42 4 "1.66"
$ws.Cells.Item(42, 5).Value = "  -6.66%  "
$ws.Cells.Item(43, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue 43 4 "0.327"
$ws.Cells.Item(43, 5).Value = "  -7.08%  "
$ws.Cells.Item(44, 5).Value = "  -6.90%  "
$ws.Cells.Item(45, 5).Value = "  -5.81%  "
Set-TextValue 46 4 "38.74"
$ws.Cells.Item(46, 5).Value = "  -2.23%  "
Set-TextValue 47 4 "141.04"
$ws.Cells.Item(47, 5).Value = "  -3.95%  "
$ws.Cells.Item(48, 5).Value = "  -6.83%  "
Set-TextValue 49 4 "0.515"
$ws.Cells.Item(49, 5).Value = "  -6.89%  "
Set-TextValue 50 4 "0.0₆0254"
$ws.Cells.Item(50, 5).Value = "  -11.78%  "
$ws.Cells.Item(51, 5).Value = "  -7.67%  "
